$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Paragraph 1: numbered "branch" list item (inherits pStyle=ListParagraph,
# numPr ilvl=0/numId=3, and rPr sz=32/szCs=32 from the preceding
# "checkout" bullet paragraph automatically).
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Range.Text = "branch – shows the current branches on the local machine if you add -a, it will show all branches even the remote branches, if you type git branch -d <branch name>, it will delete the specified branch"

# ------------------------------------------------------------------
# Paragraph 2: "git branch -m ..." — ListParagraph style, no numbering,
# ind left=1080, keeps the sz=32/szCs=32 run formatting.
# ------------------------------------------------------------------
$p1.Range.Collapse(0)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.Text = "git branch -m <oldName> <newName>: change the name of a branch"
$p2.Range.ListFormat.RemoveNumbers()
$p2.Style = "Normal"
$p2.Style = "List Paragraph"
$p2.LeftIndent = 54
$p2.Range.Font.Size = 16
$p2.Range.Font.SizeBi = 16

# ------------------------------------------------------------------
# Paragraph 3: empty ListParagraph, ind left=1080, sz=32/szCs=32 on
# the paragraph mark only (no runs).
# ------------------------------------------------------------------
$p2.Range.Collapse(0)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.Text = "x"
$p3.Range.Font.Size = 16
$p3.Range.Font.SizeBi = 16
$clearRange = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$clearRange.Text = ""

# ------------------------------------------------------------------
# Paragraph 4: same as paragraph 3 — another empty ListParagraph,
# ind left=1080.
# ------------------------------------------------------------------
$p3b = $d.Paragraphs.Last
$p3b.Range.Collapse(0)
$p3b.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last
$p4.Range.Text = "x"
$p4.Range.Font.Size = 16
$p4.Range.Font.SizeBi = 16
$clearRange2 = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$clearRange2.Text = ""

# ------------------------------------------------------------------
# Paragraph 5: ListParagraph, ind left=0, four runs of WSL setup text.
# ------------------------------------------------------------------
$p4b = $d.Paragraphs.Last
$p4b.Range.Collapse(0)
$p4b.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Last
$p5.LeftIndent = 0
$p5.Range.Font.Size = 16
$p5.Range.Font.SizeBi = 16

$run1 = "To setup git on Linux, you need to first update the package manager, otherwise it may or may not give you an error, then you will use the package manager to install git, e.g., if you use apt, it would be apt-get install git, if you are not in the root, you will need to prefix it with sudo "
$run2 = "and type in your password. Note: if you use WSL on windows, and you’ve changed the windows password for at least once after setting the WSL up, then you will probably need to type in the old windows password in order to login to the root in "
$run3 = "WSL"
$run4 = "."

$p5.Range.InsertAfter($run1)
$p5.Range.InsertAfter($run2)
$p5.Range.InsertAfter($run3)
$p5.Range.InsertAfter($run4)

Write-Host "Final paragraph count: $($d.Paragraphs.Count)"
Write-Host "Last paragraph text: $($d.Paragraphs.Last.Range.Text)"
